$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18: label for "Slope (func)" row (the existing SLOPE() based row)
$ws.Range("A18").Value = "Slope (func)"

# Row 19: replace the old "by-hand" slope estimate row with new content
$ws.Range("A19").Value = "Slope (by hand)"

$ws.Range("B19").Formula = "=(B17-B2)/(15-1)"
$ws.Range("C19").Formula = "=(C17-C2)/(15-1)"
$ws.Range("D19").Formula = "=(D17-D2)/(15-1)"
$ws.Range("E19").Formula = "=(E17-E2)/(15-1)"
$ws.Range("F19").Formula = "=(F17-F2)/(15-1)"

# Move the "stdev" label that used to sit at F19 over to H19
$ws.Range("H19").Value = "stdev"

# New summary block at H24:I25
$ws.Range("H24").Value = "expressionStability() result:"
$ws.Range("H25").Value = "old func"
$ws.Range("I25").Formula = "=STDEV(G2:G17)/(MAX(G2:G17)-MIN(G2:G17))"

# Column A width (Excel auto "best fit" width for the new longer labels)
$ws.Columns.Item(1).ColumnWidth = 13

# Recalculate so cached formula values match
$excel.Calculate()

# Restore the view: scroll so row 13 is at top, select the new summary cells
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("H24:I25").Select()
